$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Uses System"
$ws.Range("K1").Value = "Has Subprocess"
$ws.Range("L1").Value = "Related Document"
$ws.Range("M1").Value = "To Be Signed By"
$ws.Range("N1").Value = "Has Signed"
$ws.Range("O1").Value = "Approved by (filled out by PMO)"
$ws.Range("P1").Value = "Created By"
$ws.Range("Q1").Value = "Modified By"
$ws.Range("R1").Value = "Last Modified"
$ws.Range("S1").Value = "Capability_Capability"
$ws.Range("T1").Value = "MC2 Link_Capability"
$ws.Range("U1").Value = "Implements Strategy"
$ws.Range("V1").Value = "Implemented by Process"
$ws.Range("W1").Value = "Defined by Department"
$ws.Range("X1").Value = "Department"
$ws.Range("Y1").Value = "MC2 Link"
$ws.Range("Z1").Value = "Level"
$ws.Range("AA1").Value = "Belongs to Department"
$ws.Range("AB1").Value = "Has Employee"
$ws.Range("AC1").Value = "Has Leader"
$ws.Range("AD1").Value = "Belongs to Department_Class"
$ws.Range("AE1").Value = "Capability_Department"
$ws.Range("AF1").Value = "Real SOP Writer"
$ws.Range("AG1").Value = "Approver"
$ws.Range("AH1").Value = "NO Counterpart"
$ws.Range("AI1").Value = "SE Counterpart"
$ws.Range("AJ1").Value = "DK Counterpart"
$ws.Range("AK1").Value = "NL Counterpart"
$ws.Range("J2").Value = 18
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 29
$ws.Range("N2").Value = 29
$ws.Range("O2").Value = 0
$ws.Range("AK2").Value = 29
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 17
$ws.Range("M3").Value = 62
$ws.Range("N3").Value = 65
$ws.Range("O3").Value = 0
$ws.Range("V3").Value = 65
$ws.Range("AA3").Value = 62
$ws.Range("AB3").Value = 15
$ws.Range("AK3").Value = 62
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = 18
$ws.Range("N4").Value = 18
$ws.Range("O4").Value = 0
$ws.Range("AK4").Value = 18
$ws.Range("B5").Value = "EB Service Operation"
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 25
$ws.Range("O5").Value = 0
$ws.Range("W5").Value = 25
$ws.Range("Y5").Value = 10
$ws.Range("AK5").Value = 10
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = 27
$ws.Range("N6").Value = 27
$ws.Range("O6").Value = 0
$ws.Range("AK6").Value = 27
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = 24
$ws.Range("N7").Value = 24
$ws.Range("O7").Value = 0
$ws.Range("AK7").Value = 24
$ws.Range("J8").Value = 0
$ws.Range("M8").Value = 8
$ws.Range("N8").Value = 8
$ws.Range("O8").Value = 1
$ws.Range("AK8").Value = 8
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1
$ws.Range("M9").Value = 5
$ws.Range("N9").Value = 6
$ws.Range("O9").Value = 0
$ws.Range("V9").Value = 6
$ws.Range("AK9").Value = 5
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 0
$ws.Range("Z10").Value = 7
$ws.Range("AA10").Value = 0
$ws.Range("AK10").Value = 7
$ws.Range("B11").Value = "COE Product Design"
$ws.Range("J11").Value = 0
$ws.Range("M11").Value = 5
$ws.Range("N11").Value = 5
$ws.Range("O11").Value = 0
$ws.Range("Z11").Value = 5
$ws.Range("AA11").Value = 0
$ws.Range("AK11").Value = 5
$ws.Range("C12").Value = 45
$ws.Range("D12").Value = 45
$ws.Range("E12").Value = 45
$ws.Range("F12").Value = 41
$ws.Range("G12").Value = 41
$ws.Range("H12").Value = 41
$ws.Range("I12").Value = 34
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 9
$ws.Range("L12").Value = 33
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 45
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 45
$ws.Range("Q12").Value = 45
$ws.Range("R12").Value = 45
$ws.Range("S12").Value = 45
$ws.Range("T12").Value = 45
$ws.Range("U12").Value = 45
$ws.Range("V12").Value = 45
$ws.Range("W12").Value = 45
$ws.Range("X12").Value = 45
$ws.Range("Y12").Value = 45
$ws.Range("Z12").Value = 45
$ws.Range("AA12").Value = 0
$ws.Range("AB12").Value = 45
$ws.Range("AC12").Value = 45
$ws.Range("AD12").Value = 45
$ws.Range("AE12").Value = 45
$ws.Range("AF12").Value = 45
$ws.Range("AG12").Value = 45
$ws.Range("AH12").Value = 45
$ws.Range("AI12").Value = 45
$ws.Range("AJ12").Value = 45
$ws.Range("AK12").Value = 45
$ws.Range("C13").Value = 19
$ws.Range("D13").Value = 19
$ws.Range("E13").Value = 19
$ws.Range("G13").Value = 16
$ws.Range("H13").Value = 13
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 7
$ws.Range("M13").Value = 19
$ws.Range("N13").Value = 19
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 19
$ws.Range("Q13").Value = 19
$ws.Range("R13").Value = 19
$ws.Range("S13").Value = 19
$ws.Range("T13").Value = 19
$ws.Range("U13").Value = 19
$ws.Range("V13").Value = 19
$ws.Range("W13").Value = 19
$ws.Range("X13").Value = 19
$ws.Range("Y13").Value = 19
$ws.Range("Z13").Value = 19
$ws.Range("AA13").Value = 0
$ws.Range("AB13").Value = 19
$ws.Range("AC13").Value = 19
$ws.Range("AD13").Value = 19
$ws.Range("AE13").Value = 19
$ws.Range("AF13").Value = 19
$ws.Range("AG13").Value = 19
$ws.Range("AH13").Value = 19
$ws.Range("AI13").Value = 19
$ws.Range("AJ13").Value = 19
$ws.Range("AK13").Value = 19
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 7
$ws.Range("N14").Value = 7
$ws.Range("O14").Value = 0
$ws.Range("Z14").Value = 7
$ws.Range("AA14").Value = 0
$ws.Range("AK14").Value = 7
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = 1
$ws.Range("M15").Value = 17
$ws.Range("N15").Value = 17
$ws.Range("O15").Value = 0
$ws.Range("Z15").Value = 17
$ws.Range("AA15").Value = 0
$ws.Range("AK15").Value = 17
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = 11
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = 11
$ws.Range("N16").Value = 11
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 11
$ws.Range("Q16").Value = 11
$ws.Range("R16").Value = 11
$ws.Range("S16").Value = 11
$ws.Range("T16").Value = 11
$ws.Range("U16").Value = 11
$ws.Range("V16").Value = 11
$ws.Range("W16").Value = 11
$ws.Range("X16").Value = 11
$ws.Range("Y16").Value = 11
$ws.Range("Z16").Value = 11
$ws.Range("AA16").Value = 0
$ws.Range("AB16").Value = 11
$ws.Range("AC16").Value = 11
$ws.Range("AD16").Value = 11
$ws.Range("AE16").Value = 11
$ws.Range("AF16").Value = 11
$ws.Range("AG16").Value = 11
$ws.Range("AH16").Value = 11
$ws.Range("AI16").Value = 11
$ws.Range("AJ16").Value = 11
$ws.Range("AK16").Value = 11
$ws.Range("J17").Value = 5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 8
$ws.Range("N17").Value = 8
$ws.Range("O17").Value = 0
$ws.Range("Z17").Value = 8
$ws.Range("AA17").Value = 0
$ws.Range("AK17").Value = 8
